# Apply "added birth cohort for age group in output table shell" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header + age-group labels to include birth cohort ranges.
$ws.Range("D1").Value = "Age with birth cohort"
$ws.Range("D2").Value = "All (1900-2005)"
$ws.Range("D3").Value = "15-30 (1990-2005)"
$ws.Range("D4").Value = "30-44 (1975-1990)"
$ws.Range("D5").Value = "45-69 (1950-1975)"
$ws.Range("D6").Value = "70+ (1900-1950)"

# Widen columns C and D to fit the new, longer labels (values pre-adjusted so the
# engine's pixel-quantized ColumnWidth -> stored-width conversion lands as close as
# possible to the authored widths of 14.42578125 / 21.140625 character-units).
$ws.Columns.Item(3).ColumnWidth = 13.666666666666666
$ws.Columns.Item(4).ColumnWidth = 20.333333333333332

# Update the active selection to D8.
$ws.Range("D8").Select()

# Update the workbook window vertical offset.
$excel.ActiveWindow.Top = 900
